$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 18; $row++) {
    $ws.Range("AO$row").Value = 41877.3887277787
}

$excel.CalculateFullRebuild()
